$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_B2 = @"
"Fight Fire With Fire" lyrics
"@
$ws.Range("B2").Value = $v_B2

$v_C2 = @"
Do unto others as they've done to you
But what the hell is this world coming to?
Blow the universe into nothingness
Nuclear warfare shall lay us to rest
Fight fire with fire
Ending is near
Fight fire with fire
Bursting with fear
We all shall die
Time is like a fuse, short and burning fast
Armageddon's here, like said in the past
Fight fire with fire
Ending is near
Fight fire with fire
Bursting with fear
Soon to fill our lungs, the hot winds of death
The gods are laughing, so take your last breath
Fight fire with fire
Ending is near
Fight fire with fire
Bursting with fear
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight fire with fire
Fight
"@
$ws.Range("C2").Value = $v_C2

$v_D2 = @"

Metallica Lyrics

"@
$ws.Range("D2").Value = $v_D2

$v_E2 = @"
album: "Ride The Lightning" (1984)
"@
$ws.Range("E2").Value = $v_E2

$v_F2 = @"
Submit CorrectionsThanks to Grenas for correcting these lyrics.Writer(s): Lars Ulrich, James Alan Hetfield, Clifford Lee Burton
"@
$ws.Range("F2").Value = $v_F2

$v_B3 = @"
"One" lyrics
"@
$ws.Range("B3").Value = $v_B3

$v_C3 = @"
I can't remember anything
Can't tell if this is true or dream
Deep down inside I feel the scream
This terrible silence stops me
Now that the war is through with me
I'm waking up, I cannot see
That there's not much left of me
Nothing is real but pain now
Hold my breath as I wish for death
Oh please, God, wake me
Back in the womb it's much too real
In pumps life that I must feel
But can't look forward to reveal
Look to the time when I'll live
Fed through the tube that sticks in me
Just like a wartime novelty
Tied to machines that make me be
Cut this life off from me
Hold my breath as I wish for death
Oh please, God, wake me
Now the world is gone, I'm just one
Oh God, help me
Hold my breath as I wish for death
Oh please, God, help me
Darkness
Imprisoning me
All that I see
Absolute horror
I cannot live
I cannot die
Trapped in myself
Body my holding cell
Landmine
Has taken my sight
Taken my speech
Taken my hearing
Taken my arms
Taken my legs
Taken my soul
Left me with life in hell
"@
$ws.Range("C3").Value = $v_C3

$v_D3 = @"

Metallica Lyrics

"@
$ws.Range("D3").Value = $v_D3

$v_E3 = @"
album: "...And Justice For All" (1988)
"@
$ws.Range("E3").Value = $v_E3

$v_F3 = @"
Submit CorrectionsThanks to Payton for correcting these lyrics.Writer(s): Lars Ulrich, James Alan Hetfield
"@
$ws.Range("F3").Value = $v_F3

$v_B4 = @"
"Fuel" lyrics
"@
$ws.Range("B4").Value = $v_B4

$v_C4 = @"
Gimme fuel
Gimme fire
Gimme that which I desire
Ooh!
Yeah!
Turn on... I see red
Adrenaline crash and crack my head
Nitro junkie, paint me dead
And I see red
One hundred plus through black and white
War horse, warhead
Fuck 'em man, white-knuckle tight
Through black and white
On I burn
Fuel is pumping engines
Burning hard, loose and clean
And on I burn
Churning my direction
Quench my thirst with gasoline
So gimme fuel
Gimme fire
Gimme that which I desire
Ooh
Turn on beyond the bone
Swallow future, spit out home
Burn your face upon the chrome
Yeah!
Take the corner, join the crash
Headlights, head on, headlines
Another junkie lives too fast
Yeah lives way too fast, fast, fast, oohhOH!
On I burn
Fuel is pumping engines
Burning hard, loose and clean
And on I burn
Churning my direction
Quench my thirst with gasoline
So gimme fuel
Gimme fire
Gimme that which I desire
Yeah-heah
White knuckle tight!
Gimme fuel
Gimme fire
My desire
On I burn
Fuel is pumping engines
Burning hard, loose and clean
And on I burn
Churning my direction
Quench my thirst with gasoline
Gimme fuel
Gimme fire
Gimme that which I desire
Ooh
On I burn
"@
$ws.Range("C4").Value = $v_C4

$v_D4 = @"

Metallica Lyrics

"@
$ws.Range("D4").Value = $v_D4

$v_E4 = @"
album: "Reload" (1997)
"@
$ws.Range("E4").Value = $v_E4

$v_F4 = @"
Submit CorrectionsThanks to chloe for correcting these lyrics.Writer(s): Eliot Kennedy, Bryan Adams
"@
$ws.Range("F4").Value = $v_F4
